$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = "latitude"
$ws.Range("L1").Value = "longitude"
$ws.Range("L1").Select()
